# Applies the "data up to 26" update: extends the Facebook
# community-state survey sheet with four more days of data
# (23-26 Jun 2020), fills in previously-missing spot values in
# earlier rows, and corrects one existing figure in G142.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct an existing value in row 142
$ws.Range("G142").Value = 12.2405909

# Fill in previously-missing AR values for rows 134-136
$ws.Range("AR134").Value = 8.8462104
$ws.Range("AR135").Value = 10.5506085
$ws.Range("AR136").Value = 14.356492

# Row 143
$ws.Range("B143").Value = 16.4300923
$ws.Range("C143").Value = 21.1817963
$ws.Range("D143").Value = 20.2397725
$ws.Range("F143").Value = 22.149584
$ws.Range("G143").Value = 12.4251331
$ws.Range("H143").Value = 13.199867
$ws.Range("I143").Value = 11.9194381
$ws.Range("J143").Value = 14.3870968
$ws.Range("K143").Value = 13.1213255
$ws.Range("L143").Value = 15.0895932
$ws.Range("M143").Value = 16.8905516
$ws.Range("O143").Value = 9.6175166
$ws.Range("P143").Value = 19.5793223
$ws.Range("Q143").Value = 15.1881238
$ws.Range("R143").Value = 13.2566567
$ws.Range("S143").Value = 18.8726161
$ws.Range("T143").Value = 15.374564
$ws.Range("U143").Value = 15.7300681
$ws.Range("V143").Value = 21.0845323
$ws.Range("W143").Value = 12.2803712
$ws.Range("X143").Value = 13.4317781
$ws.Range("Y143").Value = 9.1019311
$ws.Range("Z143").Value = 10.6787645
$ws.Range("AA143").Value = 14.4727302
$ws.Range("AB143").Value = 15.8673754
$ws.Range("AD143").Value = 23.3105273
$ws.Range("AE143").Value = 12.0806163
$ws.Range("AF143").Value = 16.6446172
$ws.Range("AG143").Value = 16.7647578
$ws.Range("AH143").Value = 18.7169639
$ws.Range("AI143").Value = 10.926004
$ws.Range("AJ143").Value = 12.7245802
$ws.Range("AK143").Value = 14.3913285
$ws.Range("AL143").Value = 13.5518031
$ws.Range("AM143").Value = 12.2356187
$ws.Range("AN143").Value = 12.9495034
$ws.Range("AO143").Value = 19.6064668
$ws.Range("AP143").Value = 11.2915817
$ws.Range("AQ143").Value = 11.3086483
$ws.Range("AS143").Value = 11.5487716
$ws.Range("AT143").Value = 21.9870886
$ws.Range("AU143").Value = 19.3609592
$ws.Range("AV143").Value = 15.189681
$ws.Range("AW143").Value = 22.1632916
$ws.Range("AX143").Value = 19.5408345
$ws.Range("AY143").Value = 13.7730427
$ws.Range("BA143").Value = 10.388223
$ws.Range("BB143").Value = 13.1412397
$ws.Range("BC143").Value = 13.8785426
$ws.Range("BD143").Value = 13.5789521
$ws.Range("BE143").Value = 15.5949872

# Row 144
$ws.Range("B144").Value = 17.0463174
$ws.Range("C144").Value = 21.1564028
$ws.Range("D144").Value = 20.2310999
$ws.Range("F144").Value = 23.0040644
$ws.Range("G144").Value = 12.775418
$ws.Range("H144").Value = 12.9222858
$ws.Range("I144").Value = 11.4963504
$ws.Range("J144").Value = 13.4214186
$ws.Range("K144").Value = 13.5416667
$ws.Range("L144").Value = 15.8324116
$ws.Range("M144").Value = 17.2471243
$ws.Range("O144").Value = 9.894304200000001
$ws.Range("P144").Value = 20.1258308
$ws.Range("Q144").Value = 16.1170344
$ws.Range("R144").Value = 13.4225995
$ws.Range("S144").Value = 18.6094995
$ws.Range("T144").Value = 16.0758593
$ws.Range("U144").Value = 15.6757648
$ws.Range("V144").Value = 21.2745764
$ws.Range("W144").Value = 11.8854505
$ws.Range("X144").Value = 13.5907417
$ws.Range("Y144").Value = 9.627515799999999
$ws.Range("Z144").Value = 10.7857541
$ws.Range("AA144").Value = 14.4854168
$ws.Range("AB144").Value = 15.3933595
$ws.Range("AD144").Value = 24.0709086
$ws.Range("AE144").Value = 13.2498915
$ws.Range("AF144").Value = 16.8520224
$ws.Range("AG144").Value = 17.1252175
$ws.Range("AH144").Value = 18.0246439
$ws.Range("AI144").Value = 10.3153304
$ws.Range("AJ144").Value = 12.4973671
$ws.Range("AK144").Value = 14.3150808
$ws.Range("AL144").Value = 14.3441225
$ws.Range("AM144").Value = 12.1483419
$ws.Range("AN144").Value = 13.002053
$ws.Range("AO144").Value = 19.538884
$ws.Range("AP144").Value = 11.5369845
$ws.Range("AQ144").Value = 11.3603624
$ws.Range("AS144").Value = 10.9600925
$ws.Range("AT144").Value = 22.9594752
$ws.Range("AU144").Value = 19.4715984
$ws.Range("AV144").Value = 15.836489
$ws.Range("AW144").Value = 22.9719406
$ws.Range("AX144").Value = 19.9300668
$ws.Range("AY144").Value = 14.1189323
$ws.Range("BA144").Value = 10.3096462
$ws.Range("BB144").Value = 13.474175
$ws.Range("BC144").Value = 14.1179391
$ws.Range("BD144").Value = 14.2830295
$ws.Range("BE144").Value = 15.089086

# Row 145 (new date row: 23 06 2020)
$ws.Range("A145").Value = "23 06 2020"
$ws.Range("B145").Value = 17.3165138
$ws.Range("C145").Value = 22.0360745
$ws.Range("D145").Value = 20.1803961
$ws.Range("F145").Value = 24.1166012
$ws.Range("G145").Value = 13.0016802
$ws.Range("H145").Value = 13.2370663
$ws.Range("I145").Value = 11.3911743
$ws.Range("J145").Value = 15.6716418
$ws.Range("K145").Value = 13.5115725
$ws.Range("L145").Value = 16.91313
$ws.Range("M145").Value = 17.6445627
$ws.Range("O145").Value = 10.3076923
$ws.Range("P145").Value = 20.4799191
$ws.Range("Q145").Value = 16.1706377
$ws.Range("R145").Value = 13.1358945
$ws.Range("S145").Value = 18.5500694
$ws.Range("T145").Value = 15.7639909
$ws.Range("U145").Value = 15.7251752
$ws.Range("V145").Value = 22.0776219
$ws.Range("W145").Value = 11.4029664
$ws.Range("X145").Value = 13.1204554
$ws.Range("Y145").Value = 9.6080725
$ws.Range("Z145").Value = 10.848253
$ws.Range("AA145").Value = 15.1191711
$ws.Range("AB145").Value = 15.5561748
$ws.Range("AD145").Value = 24.340724
$ws.Range("AE145").Value = 13.1766286
$ws.Range("AF145").Value = 16.9959595
$ws.Range("AG145").Value = 17.3701245
$ws.Range("AH145").Value = 18.1533888
$ws.Range("AI145").Value = 10.3858078
$ws.Range("AJ145").Value = 12.6874065
$ws.Range("AK145").Value = 15.0644002
$ws.Range("AL145").Value = 15.2765066
$ws.Range("AM145").Value = 12.0733555
$ws.Range("AN145").Value = 12.9759208
$ws.Range("AO145").Value = 20.3452846
$ws.Range("AP145").Value = 11.3683295
$ws.Range("AQ145").Value = 11.5125593
$ws.Range("AS145").Value = 10.8506175
$ws.Range("AT145").Value = 23.40987
$ws.Range("AU145").Value = 20.7864045
$ws.Range("AV145").Value = 16.032284
$ws.Range("AW145").Value = 24.1767046
$ws.Range("AX145").Value = 20.2066568
$ws.Range("AY145").Value = 14.1157567
$ws.Range("BA145").Value = 10.2154649
$ws.Range("BB145").Value = 13.279306
$ws.Range("BC145").Value = 13.9073914
$ws.Range("BD145").Value = 13.5028546
$ws.Range("BE145").Value = 15.2606344

# Row 146 (new date row: 24 06 2020)
$ws.Range("A146").Value = "24 06 2020"
$ws.Range("B146").Value = 17.0936749
$ws.Range("C146").Value = 22.7740105
$ws.Range("D146").Value = 21.1436974
$ws.Range("F146").Value = 25.2981312
$ws.Range("G146").Value = 13.5824354
$ws.Range("H146").Value = 13.2388967
$ws.Range("I146").Value = 11.0983632
$ws.Range("J146").Value = 15.21181
$ws.Range("K146").Value = 14.2899761
$ws.Range("L146").Value = 17.7022691
$ws.Range("M146").Value = 18.3350568
$ws.Range("O146").Value = 10.9073959
$ws.Range("P146").Value = 20.9754309
$ws.Range("Q146").Value = 17.9560189
$ws.Range("R146").Value = 12.97533
$ws.Range("S146").Value = 18.4692339
$ws.Range("T146").Value = 16.4305433
$ws.Range("U146").Value = 15.6735795
$ws.Range("V146").Value = 22.8904417
$ws.Range("W146").Value = 11.6449151
$ws.Range("X146").Value = 13.2868916
$ws.Range("Y146").Value = 9.359681399999999
$ws.Range("Z146").Value = 10.8732211
$ws.Range("AA146").Value = 15.0747685
$ws.Range("AB146").Value = 15.5158109
$ws.Range("AD146").Value = 24.8879141
$ws.Range("AE146").Value = 13.1286613
$ws.Range("AF146").Value = 17.3084322
$ws.Range("AG146").Value = 16.5407786
$ws.Range("AH146").Value = 18.2054626
$ws.Range("AI146").Value = 10.105382
$ws.Range("AJ146").Value = 12.7694734
$ws.Range("AK146").Value = 15.3658474
$ws.Range("AL146").Value = 15.3090244
$ws.Range("AM146").Value = 12.0433725
$ws.Range("AN146").Value = 13.181903
$ws.Range("AO146").Value = 20.2608306
$ws.Range("AP146").Value = 11.670171
$ws.Range("AQ146").Value = 11.2953419
$ws.Range("AS146").Value = 11.2499169
$ws.Range("AT146").Value = 24.5334944
$ws.Range("AU146").Value = 20.9010626
$ws.Range("AV146").Value = 16.1360694
$ws.Range("AW146").Value = 25.4589453
$ws.Range("AX146").Value = 20.273411
$ws.Range("AY146").Value = 13.9455758
$ws.Range("BA146").Value = 10.5447589
$ws.Range("BB146").Value = 13.4089693
$ws.Range("BC146").Value = 13.9999813
$ws.Range("BD146").Value = 13.6277023
$ws.Range("BE146").Value = 14.4864613

# Row 147 (new date row: 25 06 2020)
$ws.Range("A147").Value = "25 06 2020"
$ws.Range("B147").Value = 16.9937206
$ws.Range("C147").Value = 23.9362063
$ws.Range("D147").Value = 21.3845834
$ws.Range("F147").Value = 26.1843834
$ws.Range("G147").Value = 13.9620312
$ws.Range("H147").Value = 12.916085
$ws.Range("I147").Value = 10.4653022
$ws.Range("J147").Value = 14.4823067
$ws.Range("K147").Value = 13.9290012
$ws.Range("L147").Value = 18.8108153
$ws.Range("M147").Value = 19.1794568
$ws.Range("O147").Value = 10.9620419
$ws.Range("P147").Value = 21.4057367
$ws.Range("Q147").Value = 18.342754
$ws.Range("R147").Value = 12.8101213
$ws.Range("S147").Value = 18.1841856
$ws.Range("T147").Value = 16.5770114
$ws.Range("U147").Value = 15.9757011
$ws.Range("V147").Value = 23.6984248
$ws.Range("W147").Value = 11.4206894
$ws.Range("X147").Value = 13.4838192
$ws.Range("Y147").Value = 10.3614833
$ws.Range("Z147").Value = 11.2238289
$ws.Range("AA147").Value = 15.2991228
$ws.Range("AB147").Value = 16.4231594
$ws.Range("AD147").Value = 24.2064887
$ws.Range("AE147").Value = 13.7672723
$ws.Range("AF147").Value = 17.4648305
$ws.Range("AG147").Value = 18.0676623
$ws.Range("AH147").Value = 18.2890104
$ws.Range("AI147").Value = 9.9858607
$ws.Range("AJ147").Value = 12.6734429
$ws.Range("AK147").Value = 15.0338561
$ws.Range("AL147").Value = 15.7220826
$ws.Range("AM147").Value = 11.9733308
$ws.Range("AN147").Value = 13.5822707
$ws.Range("AO147").Value = 20.8169737
$ws.Range("AP147").Value = 12.1027921
$ws.Range("AQ147").Value = 11.3888872
$ws.Range("AS147").Value = 11.5383538
$ws.Range("AT147").Value = 25.9422008
$ws.Range("AU147").Value = 21.1213371
$ws.Range("AV147").Value = 17.0203399
$ws.Range("AW147").Value = 26.4850715
$ws.Range("AX147").Value = 20.9381511
$ws.Range("AY147").Value = 13.9121521
$ws.Range("BA147").Value = 10.8428927
$ws.Range("BB147").Value = 13.3538857
$ws.Range("BC147").Value = 14.4330377
$ws.Range("BD147").Value = 13.9714395
$ws.Range("BE147").Value = 15.4602689

# Row 148 (new date row: 26 06 2020)
$ws.Range("A148").Value = "26 06 2020"
